$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the sprint number in the two existing day headings -------------
# ("Spint( 36) - Day 1/2 - Test Case Summary" -> "Spint( 3) - Day 1/2 - ...")
$ws.Range("B2").Value = "Spint( 3) - Day 1 - Test Case Summary"
$ws.Range("B8").Value = "Spint( 3) - Day 2 - Test Case Summary"

# --- Add a new third summary block in rows 13-16 ----------------------------
# Copy the formatting of the first block (rows 2-5) down onto the new block
# so fonts/fills/borders/number-formats match exactly.
$ws.Range("B2:C2").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)

$ws.Range("B3:C3").Copy()
$ws.Range("B14:C14").PasteSpecial(-4122)

$ws.Range("B4:C4").Copy()
$ws.Range("B15:C15").PasteSpecial(-4122)

$ws.Range("B5:C5").Copy()
$ws.Range("B16:C16").PasteSpecial(-4122)

# Fill in the new block's values
$ws.Range("B13").Value = "Spint( 3) - Day 2 - Test Case Summary"

$ws.Range("B14").Value = "Total  testcase Written"
$ws.Range("C14").Value = 123

$ws.Range("B15").Value = "Total Execution"
$ws.Range("C15").Value = 173

$ws.Range("B16").Value = "Total Review"
$ws.Range("C16").Value = 173

# Merge the new heading row, same as the other two blocks
$ws.Range("B13:C13").Merge()

# Row heights for all data rows settle at 18 in the edited workbook
foreach ($r in @(2,3,4,5,8,9,10,11,13,14,15,16)) {
    $ws.Rows.Item($r).RowHeight = 18
}

# Selection ends up on G14 in the saved workbook
[void]$ws.Range("G14").Select()
